$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("re_profiles")

# --- A) Move block2 (won data): H9:K22 -> G9:J22 (1 col left), via staging (overlaps) ---
$ws.Range("H9:K22").Cut($ws.Range("AA9"))
$ws.Range("AA9:AD22").Cut($ws.Range("G9"))

# --- B) Move block3 (wof data): O9:R22 -> L9:O22 (3 col left), via staging (overlaps) ---
$ws.Range("O9:R22").Cut($ws.Range("AA9"))
$ws.Range("AA9:AD22").Cut($ws.Range("L9"))

# --- C) Move small table header (2 rows): M2:O3 -> Q9:S10 (no overlap) ---
$ws.Range("M2:O3").Cut($ws.Range("Q9"))

# --- D) Move small table data rows (4 rows), REVERSED order: M4:O7 -> Q11:S14 reversed ---
$ws.Range("M4:O4").Cut($ws.Range("Q14"))
$ws.Range("M5:O5").Cut($ws.Range("Q13"))
$ws.Range("M6:O6").Cut($ws.Range("Q12"))
$ws.Range("M7:O7").Cut($ws.Range("Q11"))

Write-Host "done moves"

foreach ($addr in @("G9","G10","G11","H11","I11","J11","K11","L11","M11","N11","O11","P11","Q9","Q10","Q11","R11","S11","Q12","R12","S12","Q13","R13","S13","Q14","R14","S14")) {
  Write-Host "$addr :" $ws.Range($addr).Value2
}
